$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (with same bold/bordered style as existing header row)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (from AC1) onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in team record data (Wins/Losses/Ties) for every data row (2-52)
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 75  # AD
    $ws.Cells.Item($r, 31).Value = 87  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
